$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'315.50"
$ws.Range("E2").Value = "'2.42%"
$ws.Range("G2").Value = "'6"
$ws.Range("D3").Value = "'41.00"
$ws.Range("E3").Value = "'-0.11%"
$ws.Range("G3").Value = "'6"
$ws.Range("D4").Value = "'5.154"
$ws.Range("E4").Value = "'-1.54%"
$ws.Range("G4").Value = "'6"
$ws.Range("D5").Value = "'0.07596"
$ws.Range("E5").Value = "'-1.00%"
$ws.Range("G5").Value = "'6"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "'1.671"
$ws.Range("E6").Value = "'1.92%"
$ws.Range("G6").Value = "'6"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.9294"
$ws.Range("E7").Value = "'1.29%"
$ws.Range("G7").Value = "'6"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.424"
$ws.Range("E8").Value = "'-0.62%"
$ws.Range("G8").Value = "'6"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1198"
$ws.Range("E9").Value = "'-2.62%"
$ws.Range("G9").Value = "'6"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1818"
$ws.Range("E10").Value = "'-0.35%"
$ws.Range("G10").Value = "'6"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09073"
$ws.Range("E11").Value = "'-0.60%"
$ws.Range("G11").Value = "'6"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.04138"
$ws.Range("E12").Value = "'-2.92%"
$ws.Range("G12").Value = "'6"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.1053"
$ws.Range("E13").Value = "'0.01%"
$ws.Range("G13").Value = "'6"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001287"
$ws.Range("E14").Value = "'1.66%"
$ws.Range("G14").Value = "'6"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005806"
$ws.Range("E15").Value = "'0.83%"
$ws.Range("G15").Value = "'6"
$ws.Range("B16").Value = "UpBots"
$ws.Range("C16").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D16").Value = "'0.007522"
$ws.Range("E16").Value = "'0.18%"
$ws.Range("G16").Value = "'6"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.332"
$ws.Range("E17").Value = "'-0.42%"
$ws.Range("G17").Value = "'6"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'4.334"
$ws.Range("E18").Value = "'0.23%"
$ws.Range("G18").Value = "'6"
$ws.Range("D19").Value = "'0.3357"
$ws.Range("E19").Value = "'0.64%"
$ws.Range("G19").Value = "'6"
$ws.Range("D20").Value = "'7.595"
$ws.Range("E20").Value = "'3.92%"
$ws.Range("G20").Value = "'6"
$ws.Range("D21").Value = "'0.1352"
$ws.Range("E21").Value = "'-2.32%"
$ws.Range("G21").Value = "'6"
$ws.Range("D22").Value = "'0.2953"
$ws.Range("E22").Value = "'2.04%"
$ws.Range("G22").Value = "'6"
$ws.Range("D23").Value = "'0.04037"
$ws.Range("E23").Value = "'-0.83%"
$ws.Range("G23").Value = "'6"
$ws.Range("D24").Value = "'0.001275"
$ws.Range("E24").Value = "'0.95%"
$ws.Range("G24").Value = "'6"
$ws.Range("D25").Value = "'0.004049"
$ws.Range("E25").Value = "'-3.82%"
$ws.Range("G25").Value = "'6"
$ws.Range("E26").Value = "'-0.14%"
$ws.Range("G26").Value = "'6"
$ws.Range("G27").Value = "'6"
$ws.Range("G28").Value = "'6"
$ws.Range("G29").Value = "'6"
$ws.Range("G30").Value = "'6"
$ws.Range("G31").Value = "'6"
$ws.Range("G32").Value = "'6"
$ws.Range("G33").Value = "'6"
$ws.Range("G34").Value = "'6"
$ws.Range("G35").Value = "'6"
$ws.Range("G36").Value = "'6"
$ws.Range("G37").Value = "'6"
$ws.Range("D38").Value = "'0.02410"
$ws.Range("E38").Value = "'-2.10%"
$ws.Range("G38").Value = "'6"
$ws.Range("D39").Value = "'0.05150"
$ws.Range("E39").Value = "'-2.78%"
$ws.Range("G39").Value = "'6"
$ws.Range("D40").Value = "'0.007740"
$ws.Range("E40").Value = "'-1.42%"
$ws.Range("G40").Value = "'6"
$ws.Range("E41").Value = "'-1.06%"
$ws.Range("G41").Value = "'6"
$ws.Range("D42").Value = "'0.007602"
$ws.Range("E42").Value = "'13.36%"
$ws.Range("G42").Value = "'6"
$ws.Range("E43").Value = "'72.51%"
$ws.Range("G43").Value = "'6"
$ws.Range("D44").Value = "'0.008596"
$ws.Range("E44").Value = "'12.15%"
$ws.Range("G44").Value = "'6"
$ws.Range("D45").Value = "'0.3410"
$ws.Range("E45").Value = "'11.47%"
$ws.Range("G45").Value = "'6"
$ws.Range("D46").Value = "'0.00006589"
$ws.Range("E46").Value = "'-2.06%"
$ws.Range("G46").Value = "'6"
$ws.Range("G47").Value = "'6"
$ws.Range("D48").Value = "'0.2754"
$ws.Range("E48").Value = "'62.12%"
$ws.Range("G48").Value = "'6"
$ws.Range("G49").Value = "'6"
$ws.Range("G50").Value = "'6"
$ws.Range("G51").Value = "'6"
